# Add a new column 'Correction ' to the Card19 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# Fix header M1: remove trailing space from "Event "
$ws.Range("M1").Value = "Event"

# Add new header N1 "Correction " (note trailing space) with same style as other headers
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N1").Value = "Correction "

# Fill M2:M12 with "nan" (previously blank inline strings)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# Leave N2:N12 blank (new empty cells), just touch them so they materialize in the
# sheet (Excel drops truly-untouched cells on save). Toggling a no-op border style
# is enough to force the cell to persist without altering its visible formatting.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Borders.LineStyle = 0
}
